# Expand the accounting-entry template from 6 columns (date, account,
# description, debit, credit, reference) to the new 10-column layout
# (document_id, date, account_code, movement, customer_identification,
# branch_office, description, cost_center, value, observations).
#
# Column "date" (B) must stay as literal text ("2024-01-01"), NOT be
# auto-converted into an Excel date serial number -- that's the whole
# point of this edit (avoid "date serialization errors" downstream).
# Prefixing the string with a leading apostrophe forces Excel to treat
# it as literal text instead of parsing it as a date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "document_id"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "account_code"
$ws.Range("D1").Value = "movement"
$ws.Range("E1").Value = "customer_identification"
$ws.Range("F1").Value = "branch_office"
$ws.Range("G1").Value = "description"
$ws.Range("H1").Value = "cost_center"
$ws.Range("I1").Value = "value"
$ws.Range("J1").Value = "observations"

# The four new header cells (G1:J1) need the same bold/bordered header
# style already applied to A1:F1. Copy the formatting across instead of
# re-building font/border/alignment by hand.
$ws.Range("A1").Copy()
$ws.Range("G1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 2 (debit leg) ----
$ws.Range("A2").Value = 27441
$ws.Range("B2").Value = "'2024-01-01"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "'11050501"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "Debit"
$ws.Range("E2").Value = "'13832081"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "Descripción Débito"
$ws.Range("H2").Value = 235
$ws.Range("I2").Value = 119000
$ws.Range("J2").Value = "Observaciones"

# ---- Row 3 (credit leg) ----
$ws.Range("A3").Value = 27441
$ws.Range("B3").Value = "'2024-01-01"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'11100501"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "Credit"
$ws.Range("E3").Value = "'13832081"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = "Descripción Crédito"
$ws.Range("H3").Value = 235
$ws.Range("I3").Value = 119000
$ws.Range("J3").Value = "Observaciones"
